$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new Timestamp column
$ws.Range("D1").Value = "Timestamp"

# Update existing timestamp value in C2
$ws.Range("C2").Value = "2025-04-23T16:56:17.282Z"

# Add new rows 3 and 4
$ws.Range("A3").Value = "User AM melakukan klik tombol Search di Stip approval"
$ws.Range("B3").Value = "Pass"
$ws.Range("D3").Value = "2025-04-23T16:56:52.942Z"

$ws.Range("A4").Value = "User AM melakukan klik tombol Search di Stip approval"
$ws.Range("B4").Value = "Pass"
$ws.Range("D4").Value = "2025-04-23T16:56:53.051Z"
